# Automatische test-sync: 2025-06-24 19:46:50
# Append the new "Productinformatie" mail-log entry (Logs!A7:G7), roll the
# Dashboard summary count for that category (Dashboard!A6:B6), and extend the
# chart series + conditional-formatting ranges so they keep covering the
# newly added row.

$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A7:G7 -----------------------------------------------------------
$wsLogs.Range("A7").Value = "Productinformatie"
$wsLogs.Range("B7").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C7").Value = "Wat is het verschil tussen product A en product B?"
$wsLogs.Range("D7").Value = "Productinformatie"
$wsLogs.Range("E7").Value = "Beste klant,`nBedankt voor uw vraag over de verschillen tussen product A en product B. Product A is een basisversie met standaardfuncties, terwijl product B de uitgebreide versie is met extra functionaliteiten en geavanceerde opties. Als u specifieke vragen heeft over bepaalde functies of prijzen, dan hoor ik graag van u.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$wsLogs.Range("F7").Value = "2025-06-24 19:46:43"
$wsLogs.Range("G7").Value = "Ja"

# The multi-line Antwoord text triggers Excel's row auto-height; put the row
# back to the sheet's standard height so row 7 matches the look of the other
# (equally multi-line) rows above it.
$wsLogs.Rows.Item(7).RowHeight = 15

# --- Dashboard!A6:B6 --------------------------------------------------------
$wsDash.Range("A6").Value = "Productinformatie"
$wsDash.Range("B6").Value = 1

# --- Extend conditional formatting ranges to include the new row ----------
$fcsD = $wsLogs.Range("D2:D6").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D7"))
}

$fcsG = $wsLogs.Range("G2:G6").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G7"))
}

# --- Extend the Dashboard chart series to include the new category row ----
$co = $wsDash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$6,Dashboard!`$B`$2:`$B`$6,1)"
